$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 123, pushing the existing
# rows 123-137 down to 124-138 (dimension grows from R137 to R138).
$ws.Range("A123").EntireRow.Insert()

# Fill the newly-inserted row 123 with the new record's data. The leading
# descriptive columns (market id/name/region/category/quality/classification)
# are constant across this entire sheet.
$ws.Range("A123").Value = 4
$ws.Range("B123").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C123").Value = "Los Lagos"
$ws.Range("D123").Value = 44491
$ws.Range("E123").Value = 10
$ws.Range("F123").Value = 100112032
$ws.Range("G123").Value = "Zapallo italiano"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 200
$ws.Range("K123").Value = 15000
$ws.Range("L123").Value = 15000
$ws.Range("M123").Value = 15000
$ws.Range("N123").Value = "$/caja 50 unidades"
$ws.Range("O123").Value = "Región de Arica y Parinacota"
$ws.Range("P123").Value = 300
$ws.Range("Q123").Value = 50
$ws.Range("R123").Value = "Hortaliza"
